$p = $ppt.ActivePresentation

# Slide 3 ("Advanced RxSwift" day overview) - the body placeholder is shape 4
# (title, slide-number field, and the Rx logo picture come first).
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(4)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 3 is the "Day 3 - ..." line whose trailing run needs to be split.
$para3 = $tr.Paragraphs(3)

# That paragraph currently reads (as a single trailing run):
#   "Binding Track Activity (show / hide 'Loading' )"
# Split it into three runs, keeping the existing text/formatting for the
# first part and appending ", Scan Operator" after the closing parenthesis:
#   1) "Binding Track Activity (show / hide 'Loading"
#   2) "' "
#   3) "), Scan Operator"
$quoteSpace = $para3.Characters(53, 2)
$quoteSpace.Text = "’ "

$closeParen = $para3.Characters(55, 1)
$closeParen.Text = "), Scan Operator"
